# The deck's live theme part ("Integral") is swapped with the dormant
# "Office Theme" part that shipped alongside it: what PowerPoint shows as
# the presentation's Theme (color scheme used by the slide master / all
# slides) becomes the original default "Office Theme" palette.
#
# PowerPoint's ColorScheme object exposes exactly the 12 theme colour
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in this fixed
# order, writable one at a time via Colors(n).RGB. Helper below packs an
# 0xRRGGBB hex string into the R + G*256 + B*65536 long PowerPoint's RGB
# colours use (there is no RGB() builtin in this host).

function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colour scheme (the values that used
# to live in theme1.xml before the swap).
$officeThemeColors = @(
    "000000",  # 1  dk1      (Background 1 / Text 1 dark)
    "FFFFFF",  # 2  lt1      (Background 1 / Text 1 light)
    "44546A",  # 3  dk2      (Background 2 / Text 2 dark)
    "E7E6E6",  # 4  lt2      (Background 2 / Text 2 light)
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgbLong $officeThemeColors[$i - 1]
}
